$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Lgi4"
$ws.Range("C2").Value = "Adam11"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06754433333333333
$ws.Range("H2").Value = 0.202633
$ws.Range("I2").Value = 0.02266275549884949
$ws.Range("J2").Value = 0.02266275549884949
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1152986666666667
$ws.Range("N2").Value = 0.345896
$ws.Range("O2").Value = 0.1156688842087241
$ws.Range("P2").Value = 0.1156688842087241
$ws.Range("Q2").Value = 0.007787771574222222
$ws.Range("R2").Value = 0.070089944168
$ws.Range("S2").Value = 0.002621375641647045
$ws.Range("T2").Value = 0.002621375641647046

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Lgi4"
$ws.Range("C3").Value = "Adam11"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06754433333333333
$ws.Range("H3").Value = 0.202633
$ws.Range("I3").Value = 0.02266275549884949
$ws.Range("J3").Value = 0.02266275549884949
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5007510000000001
$ws.Range("N3").Value = 1.502253
$ws.Range("O3").Value = 0.5023588833325865
$ws.Range("P3").Value = 0.5023588833325865
$ws.Range("Q3").Value = 0.033822892461
$ws.Range("R3").Value = 0.304406032149
$ws.Range("S3").Value = 0.01138483654564146
$ws.Range("T3").Value = 0.01138483654564146

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lgi4"
$ws.Range("C4").Value = "Adam11"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06754433333333333
$ws.Range("H4").Value = 0.202633
$ws.Range("I4").Value = 0.02266275549884949
$ws.Range("J4").Value = 0.02266275549884949
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3807496666666667
$ws.Range("N4").Value = 1.142249
$ws.Range("O4").Value = 0.3819722324586895
$ws.Range("P4").Value = 0.3819722324586895
$ws.Range("Q4").Value = 0.02571748240188889
$ws.Range("R4").Value = 0.231457341617
$ws.Range("S4").Value = 0.008656543311560979
$ws.Range("T4").Value = 0.008656543311560981

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Lgi4"
$ws.Range("C5").Value = "Adam11"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.302381666666667
$ws.Range("H5").Value = 6.907145
$ws.Range("I5").Value = 0.7725046677002302
$ws.Range("J5").Value = 0.7725046677002302
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1152986666666667
$ws.Range("N5").Value = 0.345896
$ws.Range("O5").Value = 0.1156688842087241
$ws.Range("P5").Value = 0.1156688842087241
$ws.Range("Q5").Value = 0.2654615363244445
$ws.Range("R5").Value = 2.38915382692
$ws.Range("S5").Value = 0.08935475295891677
$ws.Range("T5").Value = 0.08935475295891677

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Lgi4"
$ws.Range("C6").Value = "Adam11"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.302381666666667
$ws.Range("H6").Value = 6.907145
$ws.Range("I6").Value = 0.7725046677002302
$ws.Range("J6").Value = 0.7725046677002302
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5007510000000001
$ws.Range("N6").Value = 1.502253
$ws.Range("O6").Value = 0.5023588833325865
$ws.Range("P6").Value = 0.5023588833325865
$ws.Range("Q6").Value = 1.152919921965
$ws.Range("R6").Value = 10.376279297685
$ws.Range("S6").Value = 0.3880745822350984
$ws.Range("T6").Value = 0.3880745822350984

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Lgi4"
$ws.Range("C7").Value = "Adam11"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.302381666666667
$ws.Range("H7").Value = 6.907145
$ws.Range("I7").Value = 0.7725046677002302
$ws.Range("J7").Value = 0.7725046677002302
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3807496666666667
$ws.Range("N7").Value = 1.142249
$ws.Range("O7").Value = 0.3819722324586895
$ws.Range("P7").Value = 0.3819722324586895
$ws.Range("Q7").Value = 0.8766310521227779
$ws.Range("R7").Value = 7.889679469105
$ws.Range("S7").Value = 0.295075332506215
$ws.Range("T7").Value = 0.295075332506215

$ws.Range("A8").Value = "ECs"
$ws.Range("B8").Value = "Lgi4"
$ws.Range("C8").Value = "Adam11"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6104853333333333
$ws.Range("H8").Value = 1.831456
$ws.Range("I8").Value = 0.2048325768009203
$ws.Range("J8").Value = 0.2048325768009203
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1152986666666667
$ws.Range("N8").Value = 0.345896
$ws.Range("O8").Value = 0.1156688842087241
$ws.Range("P8").Value = 0.1156688842087241
$ws.Range("Q8").Value = 0.07038814495288888
$ws.Range("R8").Value = 0.6334933045759999
$ws.Range("S8").Value = 0.02369275560816023
$ws.Range("T8").Value = 0.02369275560816023

$ws.Range("A9").Value = "ECs"
$ws.Range("B9").Value = "Lgi4"
$ws.Range("C9").Value = "Adam11"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6104853333333333
$ws.Range("H9").Value = 1.831456
$ws.Range("I9").Value = 0.2048325768009203
$ws.Range("J9").Value = 0.2048325768009203
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5007510000000001
$ws.Range("N9").Value = 1.502253
$ws.Range("O9").Value = 0.5023588833325865
$ws.Range("P9").Value = 0.5023588833325865
$ws.Range("Q9").Value = 0.305701141152
$ws.Range("R9").Value = 2.751310270368
$ws.Range("S9").Value = 0.1028994645518466
$ws.Range("T9").Value = 0.1028994645518466

$ws.Range("A10").Value = "ECs"
$ws.Range("B10").Value = "Lgi4"
$ws.Range("C10").Value = "Adam11"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6104853333333333
$ws.Range("H10").Value = 1.831456
$ws.Range("I10").Value = 0.2048325768009203
$ws.Range("J10").Value = 0.2048325768009203
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3807496666666667
$ws.Range("N10").Value = 1.142249
$ws.Range("O10").Value = 0.3819722324586895
$ws.Range("P10").Value = 0.3819722324586895
$ws.Range("Q10").Value = 0.2324420871715556
$ws.Range("R10").Value = 2.091978784544
$ws.Range("S10").Value = 0.0782403566409135
$ws.Range("T10").Value = 0.0782403566409135
